$d = $word.ActiveDocument

$d.Content.Find.Execute("977×2=", $true, $false, $false, $false, $false, $true, 1, $false, "558×9=", 2) | Out-Null
$d.Content.Find.Execute("309×4=", $true, $false, $false, $false, $false, $true, 1, $false, "631×6=", 2) | Out-Null
$d.Content.Find.Execute("536×7=", $true, $false, $false, $false, $false, $true, 1, $false, "216×4=", 2) | Out-Null
$d.Content.Find.Execute("359×8=", $true, $false, $false, $false, $false, $true, 1, $false, "160×4=", 2) | Out-Null
$d.Content.Find.Execute("605×7=", $true, $false, $false, $false, $false, $true, 1, $false, "538×4=", 2) | Out-Null
$d.Content.Find.Execute("606×2=", $true, $false, $false, $false, $false, $true, 1, $false, "388×7=", 2) | Out-Null
$d.Content.Find.Execute("548×4=", $true, $false, $false, $false, $false, $true, 1, $false, "274×2=", 2) | Out-Null
$d.Content.Find.Execute("130×9=", $true, $false, $false, $false, $false, $true, 1, $false, "298×7=", 2) | Out-Null
$d.Content.Find.Execute("449×2=", $true, $false, $false, $false, $false, $true, 1, $false, "417×2=", 2) | Out-Null
$d.Content.Find.Execute("239×9=", $true, $false, $false, $false, $false, $true, 1, $false, "334×6=", 2) | Out-Null
$d.Content.Find.Execute("138×5=", $true, $false, $false, $false, $false, $true, 1, $false, "394×6=", 2) | Out-Null
$d.Content.Find.Execute("143×6=", $true, $false, $false, $false, $false, $true, 1, $false, "820×2=", 2) | Out-Null
$d.Content.Find.Execute("577×8=", $true, $false, $false, $false, $false, $true, 1, $false, "217×7=", 2) | Out-Null
$d.Content.Find.Execute("341×9=", $true, $false, $false, $false, $false, $true, 1, $false, "337×9=", 2) | Out-Null
$d.Content.Find.Execute("903×6=", $true, $false, $false, $false, $false, $true, 1, $false, "959×5=", 2) | Out-Null
$d.Content.Find.Execute("928×4=", $true, $false, $false, $false, $false, $true, 1, $false, "559×6=", 2) | Out-Null
$d.Content.Find.Execute("610×4=", $true, $false, $false, $false, $false, $true, 1, $false, "728×2=", 2) | Out-Null
$d.Content.Find.Execute("313×8=", $true, $false, $false, $false, $false, $true, 1, $false, "649×3=", 2) | Out-Null
$d.Content.Find.Execute("468×2=", $true, $false, $false, $false, $false, $true, 1, $false, "289×5=", 2) | Out-Null
$d.Content.Find.Execute("592×8=", $true, $false, $false, $false, $false, $true, 1, $false, "920×9=", 2) | Out-Null
$d.Content.Find.Execute("661×9=", $true, $false, $false, $false, $false, $true, 1, $false, "650×4=", 2) | Out-Null
$d.Content.Find.Execute("955×9=", $true, $false, $false, $false, $false, $true, 1, $false, "111×8=", 2) | Out-Null
$d.Content.Find.Execute("583×7=", $true, $false, $false, $false, $false, $true, 1, $false, "686×4=", 2) | Out-Null
$d.Content.Find.Execute("466×8=", $true, $false, $false, $false, $false, $true, 1, $false, "461×8=", 2) | Out-Null
$d.Content.Find.Execute("297×2=", $true, $false, $false, $false, $false, $true, 1, $false, "995×6=", 2) | Out-Null
